# Add a new "Youtube" worksheet between "Pentester Lab" and "Udemy",
# populate it with an OWASP Top 10 2017 video link, and make it the
# active/selected sheet (mirrors the author's "Youtube section added"
# commit).

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after "Pentester Lab" (i.e. right before
# "Udemy") and rename it.
$pentesterLab = $wb.Worksheets.Item("Pentester Lab")
$youtube = $wb.Worksheets.Add($null, $pentesterLab)
$youtube.Name = "Youtube"

# Write C2 before B2 so the shared-string table picks up the URL ahead
# of the label, matching the source ordering.
$youtube.Range("C2").Value = "https://www.youtube.com/playlist?list=PLpNYlUeSK_rnsMu5S4UGtGy2HSmZdTUNl"
$youtube.Range("B2").Value = "OWASP Top 10 2017"

# Column widths (B ~44.55 chars, C ~106.44 chars).
$youtube.Columns.Item(2).ColumnWidth = 43.666666666666664
$youtube.Columns.Item(3).ColumnWidth = 105.66666666666667

# Match the author's saved selection/cursor on the new sheet.
$youtube.Range("B4").Select()

# The newly added sheet becomes the active tab (this also clears
# tabSelected on whichever sheet previously had it, e.g. "Pentester Lab").
$youtube.Activate()
